$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells for rows 218-223 with revised monthly import figures ---
$ws.Cells.Item(218,22).Value = 2858
$ws.Cells.Item(218,23).Value = 623
$ws.Cells.Item(218,25).Value = 133
$ws.Cells.Item(218,29).Value = 23
$ws.Cells.Item(219,2).Value = 6206
$ws.Cells.Item(219,3).Value = 5735
$ws.Cells.Item(219,4).Value = 1836
$ws.Cells.Item(219,14).Value = 622
$ws.Cells.Item(219,18).Value = 31
$ws.Cells.Item(219,22).Value = 3078
$ws.Cells.Item(219,23).Value = 848
$ws.Cells.Item(219,24).Value = 400
$ws.Cells.Item(219,25).Value = 171
$ws.Cells.Item(219,28).Value = 31
$ws.Cells.Item(219,29).Value = 24
$ws.Cells.Item(219,46).Value = 85
$ws.Cells.Item(220,2).Value = 7443
$ws.Cells.Item(220,3).Value = 6882
$ws.Cells.Item(220,5).Value = 921
$ws.Cells.Item(220,14).Value = 802
$ws.Cells.Item(220,19).Value = 26
$ws.Cells.Item(220,22).Value = 3607
$ws.Cells.Item(220,23).Value = 997
$ws.Cells.Item(220,24).Value = 329
$ws.Cells.Item(220,25).Value = 348
$ws.Cells.Item(220,27).Value = 112
$ws.Cells.Item(220,28).Value = 2
$ws.Cells.Item(220,29).Value = 47
$ws.Cells.Item(220,30).Value = 2610
$ws.Cells.Item(220,32).Value = 72
$ws.Cells.Item(220,35).Value = 251
$ws.Cells.Item(221,2).Value = 6477
$ws.Cells.Item(221,3).Value = 5986
$ws.Cells.Item(221,4).Value = 1906
$ws.Cells.Item(221,5).Value = 849
$ws.Cells.Item(221,6).Value = 273
$ws.Cells.Item(221,7).Value = 123
$ws.Cells.Item(221,8).Value = 160
$ws.Cells.Item(221,9).Value = 86
$ws.Cells.Item(221,10).Value = 93
$ws.Cells.Item(221,11).Value = 341
$ws.Cells.Item(221,12).Value = 155
$ws.Cells.Item(221,13).Value = 50
$ws.Cells.Item(221,14).Value = 716
$ws.Cells.Item(221,16).Value = 128
$ws.Cells.Item(221,17).Value = 46
$ws.Cells.Item(221,21).Value = 86
$ws.Cells.Item(221,22).Value = 3191
$ws.Cells.Item(221,23).Value = 814
$ws.Cells.Item(221,24).Value = 261
$ws.Cells.Item(221,25).Value = 189
$ws.Cells.Item(221,27).Value = 121
$ws.Cells.Item(221,29).Value = 42
$ws.Cells.Item(221,30).Value = 2377
$ws.Cells.Item(221,31).Value = 400
$ws.Cells.Item(221,32).Value = 62
$ws.Cells.Item(221,33).Value = 304
$ws.Cells.Item(221,34).Value = 47
$ws.Cells.Item(221,35).Value = 190
$ws.Cells.Item(221,41).Value = 63
$ws.Cells.Item(221,42).Value = 1380
$ws.Cells.Item(221,43).Value = 238
$ws.Cells.Item(221,44).Value = 34
$ws.Cells.Item(221,46).Value = 73
$ws.Cells.Item(221,49).Value = 36
$ws.Cells.Item(221,51).Value = 311
$ws.Cells.Item(221,52).Value = 79
$ws.Cells.Item(221,53).Value = 97
$ws.Cells.Item(221,54).Value = 100
$ws.Cells.Item(222,2).Value = 7251
$ws.Cells.Item(222,3).Value = 6693
$ws.Cells.Item(222,4).Value = 2166
$ws.Cells.Item(222,5).Value = 933
$ws.Cells.Item(222,6).Value = 303
$ws.Cells.Item(222,7).Value = 101
$ws.Cells.Item(222,8).Value = 184
$ws.Cells.Item(222,9).Value = 88
$ws.Cells.Item(222,10).Value = 105
$ws.Cells.Item(222,11).Value = 487
$ws.Cells.Item(222,12).Value = 220
$ws.Cells.Item(222,13).Value = 66
$ws.Cells.Item(222,14).Value = 746
$ws.Cells.Item(222,21).Value = 106
$ws.Cells.Item(222,22).Value = 3545
$ws.Cells.Item(222,23).Value = 914
$ws.Cells.Item(222,24).Value = 396
$ws.Cells.Item(222,25).Value = 199
$ws.Cells.Item(222,27).Value = 119
$ws.Cells.Item(222,29).Value = 52
$ws.Cells.Item(222,30).Value = 2631
$ws.Cells.Item(222,31).Value = 416
$ws.Cells.Item(222,35).Value = 233
$ws.Cells.Item(222,36).Value = 37
$ws.Cells.Item(222,41).Value = 82
$ws.Cells.Item(222,42).Value = 1540
$ws.Cells.Item(222,43).Value = 277
$ws.Cells.Item(222,44).Value = 44
$ws.Cells.Item(222,49).Value = 45
$ws.Cells.Item(222,51).Value = 349
$ws.Cells.Item(222,52).Value = 103
$ws.Cells.Item(222,53).Value = 95
$ws.Cells.Item(223,2).Value = 6944
$ws.Cells.Item(223,3).Value = 6381
$ws.Cells.Item(223,4).Value = 2035
$ws.Cells.Item(223,5).Value = 774
$ws.Cells.Item(223,6).Value = 250
$ws.Cells.Item(223,8).Value = 132
$ws.Cells.Item(223,10).Value = 82
$ws.Cells.Item(223,11).Value = 537
$ws.Cells.Item(223,12).Value = 245
$ws.Cells.Item(223,13).Value = 77
$ws.Cells.Item(223,14).Value = 724
$ws.Cells.Item(223,22).Value = 3543
$ws.Cells.Item(223,23).Value = 940
$ws.Cells.Item(223,25).Value = 203
$ws.Cells.Item(223,29).Value = 35
$ws.Cells.Item(223,30).Value = 2603
$ws.Cells.Item(223,33).Value = 349
$ws.Cells.Item(223,35).Value = 222
$ws.Cells.Item(223,41).Value = 75
$ws.Cells.Item(223,42).Value = 1366
$ws.Cells.Item(223,43).Value = 253
$ws.Cells.Item(223,44).Value = 30
$ws.Cells.Item(223,51).Value = 351
$ws.Cells.Item(223,52).Value = 87

# --- Add new row 224 for period 01-07-2021 ---
$ws.Cells.Item(224,1).NumberFormat = "@"
$ws.Cells.Item(224,1).Value = "01-07-2021"
$ws.Cells.Item(224,1).Style = "Normal"
$ws.Cells.Item(224,2).Value = 8055
$ws.Cells.Item(224,3).Value = 7339
$ws.Cells.Item(224,4).Value = 2493
$ws.Cells.Item(224,5).Value = 1009
$ws.Cells.Item(224,6).Value = 345
$ws.Cells.Item(224,7).Value = 119
$ws.Cells.Item(224,8).Value = 197
$ws.Cells.Item(224,9).Value = 71
$ws.Cells.Item(224,10).Value = 100
$ws.Cells.Item(224,11).Value = 632
$ws.Cells.Item(224,12).Value = 258
$ws.Cells.Item(224,13).Value = 115
$ws.Cells.Item(224,14).Value = 853
$ws.Cells.Item(224,15).Value = 195
$ws.Cells.Item(224,16).Value = 144
$ws.Cells.Item(224,17).Value = 36
$ws.Cells.Item(224,18).Value = 28
$ws.Cells.Item(224,19).Value = 48
$ws.Cells.Item(224,20).Value = 113
$ws.Cells.Item(224,21).Value = 107
$ws.Cells.Item(224,22).Value = 4071
$ws.Cells.Item(224,23).Value = 1092
$ws.Cells.Item(224,24).Value = 340
$ws.Cells.Item(224,25).Value = 329
$ws.Cells.Item(224,26).Value = 132
$ws.Cells.Item(224,27).Value = 135
$ws.Cells.Item(224,28).Value = 3
$ws.Cells.Item(224,29).Value = 48
$ws.Cells.Item(224,30).Value = 2979
$ws.Cells.Item(224,31).Value = 482
$ws.Cells.Item(224,32).Value = 115
$ws.Cells.Item(224,33).Value = 392
$ws.Cells.Item(224,34).Value = 58
$ws.Cells.Item(224,35).Value = 258
$ws.Cells.Item(224,36).Value = 40
$ws.Cells.Item(224,37).Value = 43
$ws.Cells.Item(224,38).Value = 102
$ws.Cells.Item(224,39).Value = 30
$ws.Cells.Item(224,40).Value = 60
$ws.Cells.Item(224,41).Value = 111
$ws.Cells.Item(224,42).Value = 1491
$ws.Cells.Item(224,43).Value = 224
$ws.Cells.Item(224,44).Value = 36
$ws.Cells.Item(224,45).Value = 24
$ws.Cells.Item(224,46).Value = 103
$ws.Cells.Item(224,47).Value = 94
$ws.Cells.Item(224,48).Value = 20
$ws.Cells.Item(224,49).Value = 45
$ws.Cells.Item(224,50).Value = 1
$ws.Cells.Item(224,51).Value = 391
$ws.Cells.Item(224,52).Value = 95
$ws.Cells.Item(224,53).Value = 118
$ws.Cells.Item(224,54).Value = 88
